# ZIOP_YR_FIN.xlsx update — "Doing Updates for Financials"
#
# The sheet holds yearly financial statements for ZIOP with one column per
# fiscal period (newest-first is NOT the layout — columns run oldest->newest
# left to right starting at column D; column C holds the row captions).
# This edit adds one more (newer) reporting period: a brand-new column is
# inserted at D, every existing period shifts one column to the right
# (D->E, E->F, ... J->K), and a new trailing blank column appears at L
# (mirroring the pre-existing always-blank column that used to sit at K).
# The new column D is then populated with the figures for the new period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column at D, pushing D:K -> E:L ------------------------
$ws.Columns("D").Insert(
    [Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftToRight,
    [Microsoft.Office.Interop.Excel.XlInsertFormatOrigin]::xlFormatFromLeftOrAbove)

# The freshly inserted column comes in with the "blank/default" style; Excel's
# own behaviour when inserting inside a formatted block is to inherit the
# formatting of the column that is about to sit next to it, so re-apply the
# (now-shifted) column E's formatting (number format / font / alignment) onto
# the new column D.
$ws.Range("E7:E102").Copy() | Out-Null
$ws.Range("D7:D102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

# --- 2. Fill in the new (newest) period's figures in column D --------------
# Row -> value for the new column D (matches the existing row layout; rows
# not listed here keep the blank cell the insert already produced there).
$newDValues = @{
    7  = 43465
    8  = 100
    9  = "NA"
    10 = "NA"
    12 = 34100
    13 = 0
    14 = 0
    15 = 0
    17 = 54100
    18 = -53900
    20 = 800
    21 = -52500
    22 = 0
    23 = -53100
    24 = 0
    25 = 0
    26 = -53100
    27 = 137200
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = -800
    33 = 137200
    34 = 0
    35 = 137200
    38 = 43465
    41 = 61700
    42 = 0
    43 = 1900
    44 = 0
    45 = 20700
    46 = 84300
    47 = 0
    48 = 1100
    49 = 0
    50 = 0
    51 = 0
    52 = 9700
    53 = 0
    54 = 95100
    57 = 700
    58 = 0
    59 = 8800
    60 = 9500
    61 = 0
    62 = 0
    63 = 0
    64 = 0
    65 = 0
    66 = 9500
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = -566300
    73 = 0
    74 = 0
    75 = 0
    76 = 85600
    77 = 0
    80 = 43465
    81 = 137200
    83 = 600
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = -49500
    91 = -500
    92 = 0
    93 = 0
    94 = -500
    96 = 0
    97 = 0
    98 = 0
    99 = 0
    100 = 40300
    101 = 0
    102 = -9600
}

foreach ($row in $newDValues.Keys) {
    $ws.Cells.Item($row, 4).Value = $newDValues[$row]
}
